$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.27'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.13%'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '4.35%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.104'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.58%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05586'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.05%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.472'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.86%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8171'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.22%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8400'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1328'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.13%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06986'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.45%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.02886'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.84%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09381'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.05%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001518'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.38%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0005979'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-93.89%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006127'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.05%'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '3.60%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.038'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.67%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '4.86%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03065'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-2.97%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.25%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.749'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-0.06%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04595'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-3.05%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '2.43%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001244'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.45%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004512'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-2.87%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009599'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-1.08%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '0.51%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03640'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.64%'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006189'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.70%'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1050'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.16%'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002400'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-7.75%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008850'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '6.75%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005343'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.86%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.06%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-51.59%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '20.53%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.06%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.06%'
